$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-03-06 Wednesday" "2024-03-07 Thursday"

Replace-Text "223÷4=55, 3" "842÷6=140, 2"
Replace-Text "268÷7=38, 2" "144÷6=24, 0"
Replace-Text "933÷3=311, 0" "346÷4=86, 2"
Replace-Text "130÷4=32, 2" "748÷5=149, 3"
Replace-Text "564÷6=94, 0" "542÷6=90, 2"
Replace-Text "479÷4=119, 3" "140÷3=46, 2"
Replace-Text "770÷2=385, 0" "983÷8=122, 7"
Replace-Text "226÷5=45, 1" "151÷6=25, 1"
Replace-Text "357÷7=51, 0" "922÷3=307, 1"
Replace-Text "514÷4=128, 2" "890÷4=222, 2"
Replace-Text "295÷8=36, 7" "535÷9=59, 4"
Replace-Text "178÷4=44, 2" "699÷8=87, 3"
Replace-Text "851÷8=106, 3" "506÷4=126, 2"
Replace-Text "474÷6=79, 0" "545÷6=90, 5"
Replace-Text "502÷8=62, 6" "932÷2=466, 0"
Replace-Text "404÷9=44, 8" "270÷3=90, 0"
Replace-Text "754÷8=94, 2" "896÷9=99, 5"
Replace-Text "767÷7=109, 4" "232÷2=116, 0"
Replace-Text "195÷8=24, 3" "544÷4=136, 0"
Replace-Text "608÷8=76, 0" "360÷5=72, 0"
Replace-Text "383÷4=95, 3" "980÷8=122, 4"
Replace-Text "679÷7=97, 0" "839÷9=93, 2"
Replace-Text "290÷8=36, 2" "982÷9=109, 1"
Replace-Text "148÷9=16, 4" "497÷4=124, 1"
Replace-Text "858÷7=122, 4" "276÷9=30, 6"

Write-Output "Done"
